$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arkusz1")

# Insert a new row at position 6 (shifts existing rows 6-17 down to 7-18)
$ws.Rows.Item(6).Insert()

# Fill in column A ("funkcjonalność / test") for existing rows (skip the
# newly inserted row 6 for now, to match shared-string ordering)
$ws.Range("A2").Value = "logowanie"
$ws.Range("A3").Value = "przelew zwykly"
$ws.Range("A4").Value = "przelew do ZUS"
$ws.Range("A5").Value = "przelew do US"
$ws.Range("A7").Value = "dodanie odbiorcy krajowego"
$ws.Range("A8").Value = "dodanie i edycja odbiorcy krajowego"
$ws.Range("A9").Value = "dodnie i usuniecie odbiorcy krajowego"
$ws.Range("A10").Value = "aktywacja karty"
$ws.Range("A11").Value = "dodanie automatycznej splaty całkowitej"
$ws.Range("A12").Value = "dodanie automatycznej splaty minimalnej"
$ws.Range("A13").Value = "spłata karty"
$ws.Range("A14").Value = "usunięcie automatycznej splaty minimalnej"
$ws.Range("A15").Value = "usunięcie automatycznej splaty całkowitej"
$ws.Range("A16").Value = "zasrzeżenie karty"
$ws.Range("A17").Value = "zmiana limitu karty"
$ws.Range("A18").Value = "zmiana pin karty"

# Fill in the newly inserted row 6 with data for "przelew własny"
$ws.Range("B6").Value = "test.przelewWlasny.js"
$ws.Range("C6").Value = "rachunekNadawcy"
$ws.Range("D6").Value = "rachunekOdbiorcy,tytulPrzelewu,kwota,hasloSms"
$ws.Range("A6").Value = "przelew własny"

# Update the selected cell in the sheet view
$ws.Range("A6").Select()

$wb.Save()
